$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.408.37'
$ws.Range("E2").Value = '  +0.03%  '

$ws.Range("D3").Value = '1.848.49'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.49'
$ws.Range("E5").Value = '  +0.62%  '

$ws.Range("E6").Value = '  -0.54%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07701'
$ws.Range("E8").Value = '  +2.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2919'
$ws.Range("E9").Value = '  -0.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.01'
$ws.Range("E10").Value = '  +1.84%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07744'
$ws.Range("E11").Value = '  +0.41%  '

$ws.Range("D12").Value = '1.862.22'
$ws.Range("E12").Value = '  +1.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.034'
$ws.Range("E13").Value = '  +0.64%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.00001088'
$ws.Range("E14").Value = '  +4.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6823'
$ws.Range("E15").Value = '  +0.35%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.58'
$ws.Range("E16").Value = '  +0.33%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.192'
$ws.Range("E17").Value = '  +0.34%  '

$ws.Range("D18").Value = '29.438.57'
$ws.Range("E18").Value = '  +0.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.86'
$ws.Range("E19").Value = '  +0.08%  '

$ws.Range("E20").Value = '  -0.19%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9999'
$ws.Range("E21").Value = '  -0.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.460'
$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '157.56'
$ws.Range("E24").Value = '  +0.48%  '

$ws.Range("E25").Value = '  -1.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.419'
$ws.Range("E26").Value = '  +0.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.71'
$ws.Range("E27").Value = '  +0.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.353'
$ws.Range("E28").Value = '  +4.84%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.460'
$ws.Range("E29").Value = '  +0.14%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05641'
$ws.Range("E30").Value = '  +0.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.123'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.046'
$ws.Range("E32").Value = '  +0.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.845'
$ws.Range("E33").Value = '  -0.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.164'
$ws.Range("E34").Value = '  +0.63%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7067'
$ws.Range("E35").Value = '  -0.49%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.591'
$ws.Range("E36").Value = '  +0.02%  '

$ws.Range("D37").Value = '1.226.11'
$ws.Range("E37").Value = '  -1.72%  '

$ws.Range("E38").Value = '  -1.08%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.753'
$ws.Range("E39").Value = '  -0.61%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.450'
$ws.Range("E40").Value = '  +1.21%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9066'
$ws.Range("E41").Value = '  +0.56%  '

$ws.Range("D42").Value = '2.035.82'
$ws.Range("E42").Value = '  +1.81%  '

$ws.Range("E43").Value = '  +0.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.87'
$ws.Range("E44").Value = '  +0.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.08'
$ws.Range("E45").Value = '  +0.39%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.185'
$ws.Range("E46").Value = '  +1.16%  '

$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000119'
$ws.Range("E47").Value = '  +0.13%  '

$ws.Range("E48").Value = '  +0.59%  '

$ws.Range("E49").Value = '  +3.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.029'
$ws.Range("E50").Value = '  +0.69%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.678'
$ws.Range("E51").Value = '  +0.32%  '
